$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even though it looks
# like a plain number (e.g. "-0.01"), without Excel auto-converting it to
# a numeric cell and without leaving any quote-prefix style behind.
# We do this by writing a text-returning formula and then collapsing it
# to a static value via Copy + PasteSpecial(values only) - this keeps the
# cell's format/style untouched (unlike typing a leading apostrophe, which
# stamps a quotePrefix flag onto the cell's style).
function Set-TextValue($rangeAddr, $text) {
    $ws.Range($rangeAddr).Formula = '="' + $text + '"'
    $ws.Range($rangeAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

# Column C (LF): fill top-to-bottom first, matching the order values were
# authored in (new C/A Lag, LF Lag, FFR Lag, Constant rows).
Set-TextValue "C2" "-0.159***"
Set-TextValue "C3" "-0.2***"
Set-TextValue "C4" "8.418***"
Set-TextValue "C5" "-0.007"

# Column D (FFR): fill top-to-bottom next.
Set-TextValue "D2" "-0.028***"
Set-TextValue "D3" "-0.01"
Set-TextValue "D4" "0.093"
Set-TextValue "D5" "0.212**"

# Row 6 (r2_adj): replace numeric C6, D6
$ws.Range("C6").Value = 0.85
$ws.Range("D6").Value = 0.27
